$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.822.80'
$ws.Range('E2').Value = '  +1.49%  '

$ws.Range('D3').Value = '3.767.26'
$ws.Range('E3').Value = '  -0.48%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '626.68'
$ws.Range('E5').Value = '  +3.83%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.92'
$ws.Range('E6').Value = '  +0.98%  '

$ws.Range('D7').Value = '3.764.75'
$ws.Range('E7').Value = '  -0.49%  '

$ws.Range('E8').Value = '  -0.12%  '

$ws.Range('E9').Value = '  +1.30%  '

$ws.Range('E10').Value = '  +0.62%  '

$ws.Range('E11').Value = '  +2.29%  '

$ws.Range('E12').Value = '  +0.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.05'
$ws.Range('E14').Value = '  +0.21%  '

$ws.Range('D15').Value = '4.401.52'
$ws.Range('E15').Value = '  -0.40%  '

$ws.Range('D16').Value = '3.766.99'
$ws.Range('E16').Value = '  -0.94%  '

$ws.Range('D17').Value = '68.808.61'
$ws.Range('E17').Value = '  +1.45%  '

$ws.Range('E18').Value = '  -2.99%  '

$ws.Range('E19').Value = '  -1.19%  '

$ws.Range('E20').Value = '  +0.26%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '464.99'
$ws.Range('E21').Value = '  +1.39%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.50'
$ws.Range('E22').Value = '  +0.66%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.705'
$ws.Range('E23').Value = '  +2.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.80'
$ws.Range('E24').Value = '  -0.29%  '

$ws.Range('E25').Value = '  -1.42%  '

$ws.Range('E26').Value = '  +0.88%  '

$ws.Range('E27').Value = '  +2.62%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  +0.85%  '

$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('D30').Value = '3.916.00'
$ws.Range('E30').Value = '  -0.34%  '

$ws.Range('E31').Value = '  +2.21%  '

$ws.Range('E32').Value = '  +2.45%  '

$ws.Range('E33').Value = '  -1.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.64'
$ws.Range('E34').Value = '  -1.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.169'
$ws.Range('E35').Value = '  +14.51%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.05%  '

$ws.Range('D37').Value = '3.718.11'
$ws.Range('E37').Value = '  -0.47%  '

$ws.Range('E39').Value = '  +1.99%  '

$ws.Range('E40').Value = '  +2.07%  '

$ws.Range('E41').Value = '  -0.39%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.965'
$ws.Range('E42').Value = '  -1.32%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '154.86'
$ws.Range('E45').Value = '  +1.28%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.03'
$ws.Range('E46').Value = '  -1.00%  '

$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.294'
$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.59'
$ws.Range('E48').Value = '  -1.04%  '

$ws.Range('E49').Value = '  +3.49%  '

$ws.Range('E50').Value = '  +0.51%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.37'
$ws.Range('E51').Value = '  -0.78%  '
